$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.691.02"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "2.204.30"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.35"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.09"
$ws.Range("E6").Value = "  +13.07%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.12"
$ws.Range("E10").Value = "  +7.62%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").Value = "  +6.61%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "2.534.90"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.40"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "2.201.53"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.781"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "43.646.01"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.83"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("E22").Value = "  +7.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.96"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.55"
$ws.Range("E26").Value = "  +5.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.64"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.14"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.64"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.38"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0868"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0358"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.49"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.18"
$ws.Range("E42").Value = "  +5.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.197"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.41"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.35"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("E51").Value = "  +2.42%  "
